# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 27-28) for "Vega Modelo de Temuco -
# Chirimoya", pushing the existing rows 27-36 down to rows 29-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 27; this shifts rows 27-36
# down to 29-38, carrying their formatting (incl. the date style on col D).
$ws.Rows("27:28").Insert()

# New row 27: 2021-08-30 (44438), Especial, Volumen 35, prices 3500
$ws.Cells.Item(27, 1).Value  = 10
$ws.Cells.Item(27, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(27, 3).Value  = "La Araucanía"
$ws.Cells.Item(27, 4).Value  = 44438
$ws.Cells.Item(27, 5).Value  = 9
$ws.Cells.Item(27, 6).Value  = "Fruta"
$ws.Cells.Item(27, 7).Value  = 100107
$ws.Cells.Item(27, 8).Value  = "Otros"
$ws.Cells.Item(27, 9).Value  = 100107002
$ws.Cells.Item(27, 10).Value = "Chirimoya"
$ws.Cells.Item(27, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(27, 12).Value = "Especial"
$ws.Cells.Item(27, 13).Value = 35
$ws.Cells.Item(27, 14).Value = 3500
$ws.Cells.Item(27, 15).Value = 3500
$ws.Cells.Item(27, 16).Value = 3500
$ws.Cells.Item(27, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(27, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(27, 19).Value = 3500
$ws.Cells.Item(27, 20).Value = 1

# New row 28: 2021-08-30 (44438), Primera, Volumen 20, prices 3000
$ws.Cells.Item(28, 1).Value  = 10
$ws.Cells.Item(28, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(28, 3).Value  = "La Araucanía"
$ws.Cells.Item(28, 4).Value  = 44438
$ws.Cells.Item(28, 5).Value  = 9
$ws.Cells.Item(28, 6).Value  = "Fruta"
$ws.Cells.Item(28, 7).Value  = 100107
$ws.Cells.Item(28, 8).Value  = "Otros"
$ws.Cells.Item(28, 9).Value  = 100107002
$ws.Cells.Item(28, 10).Value = "Chirimoya"
$ws.Cells.Item(28, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 20
$ws.Cells.Item(28, 14).Value = 3000
$ws.Cells.Item(28, 15).Value = 3000
$ws.Cells.Item(28, 16).Value = 3000
$ws.Cells.Item(28, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(28, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(28, 19).Value = 3000
$ws.Cells.Item(28, 20).Value = 1
